$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a numeric-looking string (e.g. barcode) into a cell while
# forcing it to be stored as TEXT (shared string) rather than being
# auto-coerced to a number, and WITHOUT leaving a new number-format style
# behind on the destination cell (it keeps the destination's original
# style). We stage the text in a scratch cell formatted as Text, copy it as
# a value into the destination (PasteSpecial values-only keeps the
# destination's existing formatting), then clear the scratch cell.
# ---------------------------------------------------------------------------
function Set-TextValue($rangeAddr, $val) {
    $scratch = $ws.Range("Z100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# Insert a new row 5 (pushing the old totals row from 5 down to 6), so we
# have room for the 4th product line. The new row inherits formatting from
# the row above it (row 4), matching the other data rows' styles.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# ---------------------------------------------------------------------------
# Row 2 — Minios PJ Mask Flakes
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("C2").Value = 44075
$ws.Range("D2").Value = 44104
$ws.Range("E2").Value = "Minios® Δημητριακα PJ Mask Flakes 250gr"
Set-TextValue "F2" "5200132750124"
$ws.Range("G2").Value = 2.95
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = "Minios"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

# ---------------------------------------------------------------------------
# Row 3 — Peppa Pig Choco Flakes
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("C3").Value = 44075
$ws.Range("D3").Value = 44104
$ws.Range("E3").Value = "Peppa Pig® Choco Flakes 250gr"
Set-TextValue "F3" "5200132750117"
$ws.Range("G3").Value = 2.95
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = "Peppa Pig"
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1.73

# ---------------------------------------------------------------------------
# Row 4 — Minios Choco Letters
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("C4").Value = 44075
$ws.Range("D4").Value = 44104
$ws.Range("E4").Value = "Minios® Δημητριακα Choco Letters 250gr"
Set-TextValue "F4" "5200132750148"
$ws.Range("G4").Value = 2.95
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = "Minios"
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 5.19

# ---------------------------------------------------------------------------
# Row 5 (new) — Servin Malaktiko Rouxon Freesh Breeze
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("C5").Value = 44075
$ws.Range("D5").Value = 44104
$ws.Range("E5").Value = "Servin® Μαλακτικο Ρουχων Freesh Breeze 750ml"
Set-TextValue "F5" "5203565995339"
$ws.Range("G5").Value = 1.75
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = "Servin"
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 2.83

# ---------------------------------------------------------------------------
# Totals row — now row 6 after the insert
# ---------------------------------------------------------------------------
$ws.Range("K6").Value = 8
$ws.Range("L6").Value = 9.75

# ---------------------------------------------------------------------------
# Extend the two color-scale conditional formats from row 4 to row 5.
# ---------------------------------------------------------------------------
$ws.Range("I1:I4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I1:I5"))
$ws.Range("J1:J4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J1:J5"))

Write-Output "done"
